$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must be forced to text
# so Excel does not silently convert them to numbers (losing formatting like trailing zeros).
$textCells = @(
    "D5",
    "D8",
    "D11",
    "D15",
    "D16",
    "D18",
    "D19",
    "D22",
    "D25",
    "D27",
    "D28",
    "D34",
    "D37",
    "D39",
    "D41",
    "D45",
    "D46",
    "D47"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "28.107.78"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "1.653.17"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "214.28"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "23.65"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "0.0876"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "1.888.14"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "1.659.22"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "0.568"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").Value = "65.87"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "28.098.48"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "233.83"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +5.25%  "
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").Value = "152.40"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.83"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Value = "1.453.78"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "0.895"
$ws.Range("E37").Value = "  +4.44%  "
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").Value = "0.935"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "69.42"
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.84"
$ws.Range("E45").Value = "  +6.48%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "5.43"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "1.796.35"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  +0.91%  "

# Restore General number format on the cells we temporarily forced to text,
# so only the value (not the display format) differs from before.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}

Write-Host "Applied cryptos list update"
